$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing existing rows 5-11 down to 6-12.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly data point.
$ws.Range("A5").Value = 9
$ws.Range("B5").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value2 = 44533
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101008
$ws.Range("J5").Value = "Mora"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 4000
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 4000
$ws.Range("Q5").Value = "$/bandeja 2 kilos"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 2000
$ws.Range("T5").Value = 2
